# Fruta / hortaliza, semanal
# Two new weekly price records (date 2021-10-05, serial 44474) are added
# for "Naranja" at "Vega Central Mapocho de Santiago". They are inserted
# right before the existing row 397, which pushes all subsequent data
# rows (old 397-461) down by two positions (to 399-463).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 397 (shifts old row 397.. down by 2)
$ws.Rows.Item(397).Insert()
$ws.Rows.Item(397).Insert()

# --- New row 397 ---
$ws.Range("A397").Value() = 9
$ws.Range("B397").Value() = "Vega Central Mapocho de Santiago"
$ws.Range("C397").Value() = "Metropolitana"
$ws.Range("D397").Value() = 44474
$ws.Range("E397").Value() = 13
$ws.Range("F397").Value() = "Fruta"
$ws.Range("G397").Value() = 100102
$ws.Range("H397").Value() = "Cítricos"
$ws.Range("I397").Value() = 100102005
$ws.Range("J397").Value() = "Naranja"
$ws.Range("K397").Value() = "Lane Late"
$ws.Range("L397").Value() = "Primera"
$ws.Range("M397").Value() = 480
$ws.Range("N397").Value() = 6000
$ws.Range("O397").Value() = 6000
$ws.Range("P397").Value() = 6000
$ws.Range("Q397").Value() = "$/malla 18 kilos"
$ws.Range("R397").Value() = "Región de O'Higgins"
$ws.Range("S397").Value() = 333
$ws.Range("T397").Value() = 18

# --- New row 398 ---
$ws.Range("A398").Value() = 9
$ws.Range("B398").Value() = "Vega Central Mapocho de Santiago"
$ws.Range("C398").Value() = "Metropolitana"
$ws.Range("D398").Value() = 44474
$ws.Range("E398").Value() = 13
$ws.Range("F398").Value() = "Fruta"
$ws.Range("G398").Value() = 100102
$ws.Range("H398").Value() = "Cítricos"
$ws.Range("I398").Value() = 100102005
$ws.Range("J398").Value() = "Naranja"
$ws.Range("K398").Value() = "Navel Late"
$ws.Range("L398").Value() = "Primera"
$ws.Range("M398").Value() = 880
$ws.Range("N398").Value() = 6500
$ws.Range("O398").Value() = 7000
$ws.Range("P398").Value() = 6670
$ws.Range("Q398").Value() = "$/caja 18 kilos granel"
$ws.Range("R398").Value() = "Provincia de Melipilla"
$ws.Range("S398").Value() = 371
$ws.Range("T398").Value() = 18
